$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 66
$ws.Range("F3").Value = 21661
$ws.Range("F5").Value = 339
$ws.Range("F6").Value = 1137
$ws.Range("F7").Value = 32
$ws.Range("F8").Value = 8063
$ws.Range("F9").Value = 565
$ws.Range("F10").Value = 50
$ws.Range("F11").Value = 772
$ws.Range("F13").Value = 72
$ws.Range("F14").Value = 195
$ws.Range("F15").Value = 189
$ws.Range("F16").Value = 37
$ws.Range("F17").Value = 238
$ws.Range("F18").Value = 13
$ws.Range("F19").Value = 1365
$ws.Range("F20").Value = 568
$ws.Range("F22").Value = 713
$ws.Range("F24").Value = 96
$ws.Range("F27").Value = 1213
$ws.Range("F28").Value = 71
$ws.Range("F31").Value = 611
$ws.Range("F33").Value = 160
$ws.Range("F34").Value = 5164
$ws.Range("F39").Value = 13318
$ws.Range("F40").Value = 1376
$ws.Range("F42").Value = 61
$ws.Range("F44").Value = 335
$ws.Range("F45").Value = 464
$ws.Range("F46").Value = 4079
$ws.Range("F47").Value = 35
$ws.Range("F48").Value = 334

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 66
$ws.Range("F3").Value = 21661
$ws.Range("F4").Value = 1137
$ws.Range("F5").Value = 32
$ws.Range("F6").Value = 8063
$ws.Range("F7").Value = 565
$ws.Range("F8").Value = 50
$ws.Range("F9").Value = 772
$ws.Range("F11").Value = 72
$ws.Range("F12").Value = 195
$ws.Range("F13").Value = 189
$ws.Range("F14").Value = 37
$ws.Range("F15").Value = 238
$ws.Range("F16").Value = 13
$ws.Range("F17").Value = 1365
$ws.Range("F18").Value = 568
$ws.Range("F20").Value = 713
$ws.Range("F22").Value = 96
$ws.Range("F25").Value = 1213
$ws.Range("F26").Value = 71
$ws.Range("F30").Value = 611
$ws.Range("F33").Value = 160
$ws.Range("F35").Value = 5164
$ws.Range("F40").Value = 13318
$ws.Range("F41").Value = 1376
$ws.Range("F42").Value = 61
$ws.Range("F44").Value = 335
$ws.Range("F45").Value = 464
$ws.Range("F46").Value = 4079
$ws.Range("F47").Value = 35
$ws.Range("F48").Value = 334
